$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text-preserving number format for numeric-looking price strings
# so Excel does not auto-convert them to actual numbers (column D values
# are stored as text in the source data, e.g. "613.17", "1.00").
$textCells = @("D5", "D6", "D11", "D13", "D14", "D20", "D21", "D22", "D23", "D24", "D26", "D28", "D29", "D30", "D33", "D34", "D35", "D41", "D45", "D46", "D47", "D48", "D49")
foreach ($cell in $textCells) {
    $ws.Range($cell).NumberFormat = "@"
}

$ws.Range("D2").Value = "67.608.35"
$ws.Range("E2").Value = "  -2.70%  "
$ws.Range("D3").Value = "3.536.78"
$ws.Range("E3").Value = "  -3.67%  "
$ws.Range("D5").Value = "613.17"
$ws.Range("E5").Value = "  -5.45%  "
$ws.Range("D6").Value = "154.17"
$ws.Range("E6").Value = "  -3.47%  "
$ws.Range("D7").Value = "3.533.29"
$ws.Range("E7").Value = "  -3.64%  "
$ws.Range("E8").Value = "  +0.00%  "
$ws.Range("E9").Value = "  -1.98%  "
$ws.Range("E10").Value = "  -2.61%  "
$ws.Range("D11").Value = "6.89"
$ws.Range("E11").Value = "  -3.09%  "
$ws.Range("E12").Value = "  -1.71%  "
$ws.Range("D13").Value = "0.0000223"
$ws.Range("E13").Value = "  -3.54%  "
$ws.Range("D14").Value = "32.23"
$ws.Range("E14").Value = "  -0.99%  "
$ws.Range("D15").Value = "4.131.59"
$ws.Range("E15").Value = "  -3.74%  "
$ws.Range("D16").Value = "3.533.02"
$ws.Range("E16").Value = "  -3.38%  "
$ws.Range("D17").Value = "67.563.99"
$ws.Range("E17").Value = "  -2.71%  "
$ws.Range("E18").Value = "  +0.58%  "
$ws.Range("E19").Value = "  -0.83%  "
$ws.Range("D20").Value = "15.60"
$ws.Range("E20").Value = "  -2.39%  "
$ws.Range("D21").Value = "454.32"
$ws.Range("E21").Value = "  -2.24%  "
$ws.Range("D22").Value = "9.43"
$ws.Range("E22").Value = "  -3.29%  "
$ws.Range("D23").Value = "0.644"
$ws.Range("E23").Value = "  +0.03%  "
$ws.Range("D24").Value = "78.74"
$ws.Range("E24").Value = "  -0.90%  "
$ws.Range("B25").Value = "WrappedeETH"
$ws.Range("C25").Value = "https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth"
$ws.Range("D25").Value = "3.671.79"
$ws.Range("E25").Value = "  -3.85%  "
$ws.Range("B26").Value = "Dai"
$ws.Range("C26").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D26").Value = "1.00"
$ws.Range("E26").Value = "  +0.08%  "
$ws.Range("E27").Value = "  -3.63%  "
$ws.Range("D28").Value = "10.50"
$ws.Range("E28").Value = "  -2.55%  "
$ws.Range("D29").Value = "8.37"
$ws.Range("E29").Value = "  -6.06%  "
$ws.Range("D30").Value = "1.70"
$ws.Range("E30").Value = "  +1.86%  "
$ws.Range("E31").Value = "  -1.59%  "
$ws.Range("B33").Value = "ImmutableX"
$ws.Range("C33").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D33").Value = "1.92"
$ws.Range("E33").Value = "  -3.89%  "
$ws.Range("B34").Value = "EthereumClassic"
$ws.Range("C34").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D34").Value = "26.03"
$ws.Range("E34").Value = "  -2.15%  "
$ws.Range("D35").Value = "6.23"
$ws.Range("E35").Value = "  -3.50%  "
$ws.Range("D37").Value = "3.530.61"
$ws.Range("E37").Value = "  -3.61%  "
$ws.Range("E38").Value = "  -4.05%  "
$ws.Range("E40").Value = "  -0.04%  "
$ws.Range("D41").Value = "171.56"
$ws.Range("E41").Value = "  -3.49%  "
$ws.Range("E42").Value = "  -4.69%  "
$ws.Range("E43").Value = "  -1.77%  "
$ws.Range("E44").Value = "  -2.86%  "
$ws.Range("D45").Value = "0.893"
$ws.Range("E45").Value = "  -3.65%  "
$ws.Range("D46").Value = "29.54"
$ws.Range("E46").Value = "  +9.48%  "
$ws.Range("D47").Value = "45.77"
$ws.Range("E47").Value = "  -1.76%  "
$ws.Range("D48").Value = "2.65"
$ws.Range("E48").Value = "  -2.58%  "
$ws.Range("D49").Value = "1.24"
$ws.Range("E49").Value = "  -2.08%  "
$ws.Range("E50").Value = "  -1.93%  "
$ws.Range("E51").Value = "  -2.72%  "
